$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value into a cell while guaranteeing it stays
# text (Excel would otherwise auto-coerce purely-numeric-looking strings,
# e.g. "233.63", into a numeric value). We briefly force the Text number
# format, assign the value, then restore the default "Normal" style so the
# cell ends up exactly like the other unstyled text cells on the sheet.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "37.780.30"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "2.083.41"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue $ws.Range("D5") "233.63"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("E6").Value = "  +0.13%  "
Set-TextValue $ws.Range("D7") "58.88"
$ws.Range("E7").Value = "  +3.17%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +2.93%  "
$ws.Range("D12").Value = "2.390.18"
$ws.Range("E12").Value = "  +0.54%  "
Set-TextValue $ws.Range("D13") "14.73"
$ws.Range("E13").Value = "  +1.71%  "
Set-TextValue $ws.Range("D14") "21.24"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("E15").Value = "  +2.44%  "
Set-TextValue $ws.Range("D16") "5.34"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "2.067.78"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "37.740.38"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  -0.24%  "
Set-TextValue $ws.Range("D20") "71.74"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("E21").Value = "  +3.20%  "
Set-TextValue $ws.Range("D22") "228.44"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("E25").Value = "  +1.43%  "
Set-TextValue $ws.Range("D26") "171.23"
$ws.Range("E26").Value = "  +0.80%  "
Set-TextValue $ws.Range("D27") "9.49"
$ws.Range("E27").Value = "  +6.57%  "
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  +2.46%  "
$ws.Range("E32").Value = "  +2.34%  "
Set-TextValue $ws.Range("D33") "0.0636"
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  +0.19%  "
Set-TextValue $ws.Range("D40") "0.0981"
$ws.Range("E40").Value = "  -0.94%  "
Set-TextValue $ws.Range("D41") "99.35"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("E42").Value = "  +2.77%  "
Set-TextValue $ws.Range("D43") "17.16"
$ws.Range("E43").Value = "  +9.83%  "
Set-TextValue $ws.Range("D44") "2.91"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").Value = "1.452.60"
$ws.Range("E45").Value = "  -0.28%  "
Set-TextValue $ws.Range("D46") "1.16"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  -4.32%  "
Set-TextValue $ws.Range("D48") "1.07"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").Value = "2.275.74"
$ws.Range("E51").Value = "  +0.32%  "
